$d = $word.ActiveDocument

# Locate the date line "Ngày 6 tháng 11 năm 2024" and collapse the
# search range onto it.
$rng = $d.Content
$rng.Find.Execute("Ngày 6 tháng 11 năm 2024", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $rng.Start
$end = $rng.End

# "Ngày 6 tháng 11 năm 2024"
#  0123456789...
# "Ngày " -> 5 chars, then the day number "6" -> 1 char, then the rest.
$dayRange = $d.Range($start + 5, $start + 6)

# Replace the day number itself.
$dayRange.Text = "15"

# Force the run containing the day number to materialize its own
# run properties (splitting it off from the neighbouring runs) by
# toggling a character property on/off.
$dayRange.Bold = 1
$dayRange.Bold = 0
